# "double induction exp. data and analysis"
#
# Adds a 5th plate-layout sheet ("h1_specificity_09.28.22") for the
# 09.28.22 double-induction (malathion +/- permethrin/zeta-cypermethrin)
# plate, appended after the existing four sheets, and leaves it as the
# active/selected tab.

$wb = $excel.ActiveWorkbook

# --- sheet that was active before (neo2_specificity_09.09.22) moves its
#     selection before focus leaves it -----------------------------------
$wsNeo2_0909 = $wb.Worksheets.Item(2)
$wsNeo2_0909.Activate()
$wsNeo2_0909.Range("I7").Select() | Out-Null

# --- neo2_specificity_09.13.22 also gets visited / its cursor reset ------
$wsNeo2_0913 = $wb.Worksheets.Item(4)
$wsNeo2_0913.Activate()
$wsNeo2_0913.Range("A1").Select() | Out-Null

# --- append the new sheet at the end -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "h1_specificity_09.28.22"

# Row 1: headers for the 09.28.22 double-induction conditions (malathion,
# malathion+permethrin, malathion+zeta-cypermethrin), strain 4736
$ws.Range("A1").Value = "4736_mal_0"
$ws.Range("B1").Value = "4736_mal_1"
$ws.Range("C1").Value = "4736_mal_2"
$ws.Range("D1").Value = "4736_mal+perm_0"
$ws.Range("E1").Value = "4736_mal+perm_1"
$ws.Range("F1").Value = "4736_mal+perm_2"
$ws.Range("G1").Value = "4736_mal+zeta_0"
$ws.Range("H1").Value = "4736_mal+zeta_1"
$ws.Range("I1").Value = "4736_mal+zeta_2"

# Row 2: existing single-induction controls/perm/zeta for strain 4736
$ws.Range("A2").Value = "4736_control_0"
$ws.Range("B2").Value = "4736_control_1"
$ws.Range("C2").Value = "4736_control_2"
$ws.Range("D2").Value = "4736_perm_0"
$ws.Range("E2").Value = "4736_perm_1"
$ws.Range("F2").Value = "4736_perm_2"
$ws.Range("G2").Value = "4736_zeta_0"
$ws.Range("H2").Value = "4736_zeta_1"
$ws.Range("I2").Value = "4736_zeta_2"

# Row 3: double-induction conditions for strain 0953
$ws.Range("A3").Value = "0953_mal_0"
$ws.Range("B3").Value = "0953_mal_1"
$ws.Range("C3").Value = "0953_mal_2"
$ws.Range("D3").Value = "0953_mal+perm_0"
$ws.Range("E3").Value = "0953_mal+perm_1"
$ws.Range("F3").Value = "0953_mal+perm_2"
$ws.Range("G3").Value = "0953_mal+zeta_0"
$ws.Range("H3").Value = "0953_mal+zeta_1"
$ws.Range("I3").Value = "0953_mal+zeta_2"

# Row 4: existing single-induction controls/perm/zeta for strain 0953
$ws.Range("A4").Value = "0953_control_0"
$ws.Range("B4").Value = "0953_control_1"
$ws.Range("C4").Value = "0953_control_2"
$ws.Range("D4").Value = "0953_perm_0"
$ws.Range("E4").Value = "0953_perm_1"
$ws.Range("F4").Value = "0953_perm_2"
$ws.Range("G4").Value = "0953_zeta_0"
$ws.Range("H4").Value = "0953_zeta_1"
$ws.Range("I4").Value = "0953_zeta_2"

# Make the new sheet the active/visible tab, with its own cursor position
# and zoom level (170%), matching the plate-reader viewing preset.
$ws.Activate()
$ws.Range("F18").Select() | Out-Null
$excel.ActiveWindow.Zoom = 170
